$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D column values that look like plain numbers need a leading apostrophe
# so Excel stores them as text (matching the source inline-string data),
# exactly like typing e.g. `'690.99` into a cell.

$ws.Range("D2").Value = "71.327.17"
$ws.Range("E2").Value = "  +1.52%  "
$ws.Range("D3").Value = "3.872.82"
$ws.Range("E3").Value = "  +2.07%  "
$ws.Range("D5").Value = "'690.99"
$ws.Range("E5").Value = "  +3.02%  "
$ws.Range("D6").Value = "'173.23"
$ws.Range("E6").Value = "  +2.65%  "
$ws.Range("D7").Value = "3.871.97"
$ws.Range("E7").Value = "  +2.06%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").Value = "'0.164"
$ws.Range("E10").Value = "  +1.84%  "
$ws.Range("D11").Value = "'7.41"
$ws.Range("E11").Value = "  +4.26%  "
$ws.Range("D12").Value = "'0.463"
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("E13").Value = "  +6.34%  "
$ws.Range("D14").Value = "'36.70"
$ws.Range("E14").Value = "  +2.91%  "
$ws.Range("D15").Value = "4.522.81"
$ws.Range("E15").Value = "  +2.05%  "
$ws.Range("D16").Value = "3.868.49"
$ws.Range("E16").Value = "  +2.25%  "
$ws.Range("D17").Value = "71.346.88"
$ws.Range("E17").Value = "  +1.50%  "
$ws.Range("D18").Value = "'17.84"
$ws.Range("E18").Value = "  +1.10%  "
$ws.Range("D19").Value = "'7.27"
$ws.Range("E19").Value = "  +1.14%  "
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("D21").Value = "'11.08"
$ws.Range("E21").Value = "  -3.74%  "
$ws.Range("D22").Value = "'492.91"
$ws.Range("E22").Value = "  +3.79%  "
$ws.Range("D23").Value = "'0.724"
$ws.Range("E23").Value = "  +1.49%  "
$ws.Range("D24").Value = "'84.94"
$ws.Range("E24").Value = "  +2.42%  "
$ws.Range("E25").Value = "  +4.31%  "
$ws.Range("D26").Value = "'12.40"
$ws.Range("E26").Value = "  +1.55%  "
$ws.Range("D27").Value = "'10.57"
$ws.Range("E27").Value = "  +2.28%  "
$ws.Range("D28").Value = "'2.16"
$ws.Range("E28").Value = "  +1.82%  "
$ws.Range("D29").Value = "4.024.63"
$ws.Range("E29").Value = "  +2.06%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("E31").Value = "  +9.56%  "
$ws.Range("E32").Value = "  +3.07%  "
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("D34").Value = "'29.91"
$ws.Range("E34").Value = "  +1.17%  "
$ws.Range("E35").Value = "  +0.54%  "
$ws.Range("D36").Value = "'9.34"
$ws.Range("E36").Value = "  +2.56%  "
$ws.Range("D37").Value = "3.823.71"
$ws.Range("E37").Value = "  +1.89%  "
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("E39").Value = "  +2.39%  "
$ws.Range("E40").Value = "  +13.03%  "
$ws.Range("E41").Value = "  +2.10%  "
$ws.Range("D42").Value = "'6.09"
$ws.Range("E42").Value = "  +2.13%  "
$ws.Range("D43").Value = "'1.03"
$ws.Range("E43").Value = "  +6.62%  "
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D46").Value = "'163.87"
$ws.Range("E46").Value = "  +3.03%  "
$ws.Range("D47").Value = "'0.000310"
$ws.Range("E47").Value = "  +7.79%  "
$ws.Range("E48").Value = "  +1.41%  "
$ws.Range("D49").Value = "'44.59"
$ws.Range("E49").Value = "  -1.96%  "
$ws.Range("E50").Value = "  +1.46%  "
$ws.Range("E51").Value = "  -2.01%  "
